$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-table refresh (coinranking.com snapshot): each row below carries the
# updated Price (D) and/or Volume(1h) (E) text for that coin, plus two rows
# (36/37 and 40/41) whose Coin/Link/Price/Volume data got reordered/swapped.
# D-column values that look like plain numbers ("205.71", "1.00", "0.0853", ...)
# are written with a leading apostrophe so Excel keeps them as literal TEXT
# (matching the source data, e.g. "1.00" must stay "1.00", not become 1) —
# values such as "26.707.31" already fail numeric parsing (two dots) and need no help.

# Row 2
$ws.Range("D2").Value = '26.707.31'
$ws.Range("E2").Value = '  +0.11%  '

# Row 3
$ws.Range("D3").Value = '1.533.74'
$ws.Range("E3").Value = '  -1.48%  '

# Row 4
$ws.Range("E4").Value = '  -0.15%  '

# Row 5
$ws.Range("D5").Value = '''205.71'
$ws.Range("E5").Value = '  -0.06%  '

# Row 6
$ws.Range("E6").Value = '  -1.08%  '

# Row 7
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("D8").Value = '''21.35'
$ws.Range("E8").Value = '  -2.66%  '

# Row 9
$ws.Range("E9").Value = '  -1.21%  '

# Row 10
$ws.Range("E10").Value = '  -0.51%  '

# Row 11
$ws.Range("D11").Value = '''0.0853'
$ws.Range("E11").Value = '  -0.97%  '

# Row 12
$ws.Range("E12").Value = '  -1.44%  '

# Row 13
$ws.Range("D13").Value = '1.526.39'
$ws.Range("E13").Value = '  -2.26%  '

# Row 14
$ws.Range("E14").Value = '  -1.73%  '

# Row 15
$ws.Range("E15").Value = '  -1.04%  '

# Row 16
$ws.Range("D16").Value = '26.709.82'

# Row 17
$ws.Range("D17").Value = '''61.23'
$ws.Range("E17").Value = '  -0.68%  '

# Row 18
$ws.Range("D18").Value = '''212.11'
$ws.Range("E18").Value = '  -0.60%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0681'
$ws.Range("E19").Value = '  +1.03%  '

# Row 20
$ws.Range("D20").Value = '''7.20'
$ws.Range("E20").Value = '  -1.79%  '

# Row 21
$ws.Range("E21").Value = '  -0.19%  '

# Row 22
$ws.Range("E22").Value = '  -1.87%  '

# Row 23
$ws.Range("E23").Value = '  -2.87%  '

# Row 24
$ws.Range("E24").Value = '  -2.60%  '

# Row 25
$ws.Range("D25").Value = '''151.73'
$ws.Range("E25").Value = '  -0.31%  '

# Row 26
$ws.Range("E26").Value = '  -3.24%  '

# Row 27
$ws.Range("D27").Value = '''14.83'
$ws.Range("E27").Value = '  +0.13%  '

# Row 28
$ws.Range("E28").Value = '  -0.15%  '

# Row 29
$ws.Range("E29").Value = '  -1.22%  '

# Row 30
$ws.Range("E30").Value = '  -1.04%  '

# Row 31
$ws.Range("E31").Value = '  -2.08%  '

# Row 32
$ws.Range("E32").Value = '  +2.83%  '

# Row 33
$ws.Range("D33").Value = '1.359.89'
$ws.Range("E33").Value = '  -1.84%  '

# Row 34
$ws.Range("E34").Value = '  +0.02%  '

# Row 35
$ws.Range("E35").Value = '  -3.35%  '

# Row 36
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '''2.27'
$ws.Range("E36").Value = '  -0.53%  '

# Row 37
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '''0.937'
$ws.Range("E37").Value = '  +0.72%  '

# Row 38
$ws.Range("E38").Value = '  +0.11%  '

# Row 39
$ws.Range("E39").Value = '  +0.90%  '

# Row 40
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = '''1.00'
$ws.Range("E40").Value = '  -0.15%  '

# Row 41
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = '''0.797'
$ws.Range("E41").Value = '  -1.70%  '

# Row 42
$ws.Range("D42").Value = '''5.66'
$ws.Range("E42").Value = '  +5.22%  '

# Row 43
$ws.Range("D43").Value = '''0.992'
$ws.Range("E43").Value = '  -0.20%  '

# Row 44
$ws.Range("D44").Value = '''2.19'
$ws.Range("E44").Value = '  +0.48%  '

# Row 45
$ws.Range("E45").Value = '  -1.02%  '

# Row 46
$ws.Range("D46").Value = '''62.51'
$ws.Range("E46").Value = '  -0.86%  '

# Row 47
$ws.Range("D47").Value = '1.666.67'

# Row 48
$ws.Range("D48").Value = '''85.31'
$ws.Range("E48").Value = '  -0.06%  '

# Row 49
$ws.Range("D49").Value = '''0.0505'
$ws.Range("E49").Value = '  +2.36%  '

# Row 50
$ws.Range("D50").Value = '0.0₇0969'
$ws.Range("E50").Value = '  -0.39%  '
